{"js": "// Load all paragraphs in the document body so we can locate the\n// \"1 Introduction\" heading and the following introduction paragraph.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items,text,style\");\nawait context.sync();\n\n// Find the Heading 1 paragraph whose text is exactly \"1 Introduction\" and\n// the paragraph right after it (the \"In hopes to provide safety data...\"\n// introduction text).\nlet headingPara = null;\nlet introPara = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const p = paragraphs.items[i];\n  if (headingPara === null && p.style === \"Heading 1\" && p.text.trim() === \"1 Introduction\") {\n    headingPara = p;\n    introPara = paragraphs.items[i + 1];\n    break;\n  }\n}\n\nif (headingPara === null || introPara === null) {\n  throw new Error(\"Could not locate the '1 Introduction' heading paragraph.\");\n}\n\n// Prepend the raw placeholder text (\"\\newpage #1 Introduction \") to the\n// introduction paragraph as four separate runs, matching the way the\n// existing document builds up a paragraph's text from multiple runs.\nconst ooxml =\n  '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  \"<pkg:xmlData>\" +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  \"<w:body>\" +\n  \"<w:p>\" +\n  '<w:r><w:t xml:space=\"preserve\">\\\\newpage</w:t></w:r>' +\n  '<w:r><w:t xml:space=\"preserve\"> </w:t></w:r>' +\n  '<w:r><w:t xml:space=\"preserve\">#1 Introduction</w:t></w:r>' +\n  '<w:r><w:t xml:space=\"preserve\"> </w:t></w:r>' +\n  \"</w:p>\" +\n  \"</w:body>\" +\n  \"</w:document>\" +\n  \"</pkg:xmlData>\" +\n  \"</pkg:part>\" +\n  \"</pkg:package>\";\n\nintroPara.insertOoxml(ooxml, Word.InsertLocation.start);\n\n// Remove the standalone \"1 Introduction\" heading paragraph (and its\n// bookmark) entirely -- the heading text now lives, unformatted, at the\n// start of the paragraph above.\nheadingPara.delete();\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Locate the \"1 Introduction\" Heading 1 paragraph, and the paragraph right\n# after it (the \"In hopes to provide safety data...\" introduction text).\n$headingPara = $null\n$introPara = $null\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    if ($p.Style.NameLocal -eq \"Heading 1\" -and $p.Range.Text.Trim() -eq \"1 Introduction\") {\n        $headingPara = $p\n        $introPara = $d.Paragraphs.Item($i + 1)\n        break\n    }\n}\n\nif ($headingPara -eq $null -or $introPara -eq $null) {\n    throw \"Could not locate the '1 Introduction' heading paragraph.\"\n}\n\n# Prepend the raw placeholder text (\"\\newpage #1 Introduction \") to the\n# introduction paragraph as four separate runs, matching the way the rest\n# of the document builds up a paragraph's text from multiple runs.\n$insertPoint = $d.Range($introPara.Range.Start, $introPara.Range.Start)\n$xml = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?><pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body><w:p><w:r><w:t xml:space=\"preserve\">\\newpage</w:t></w:r><w:r><w:t xml:space=\"preserve\"> </w:t></w:r><w:r><w:t xml:space=\"preserve\">#1 Introduction</w:t></w:r><w:r><w:t xml:space=\"preserve\"> </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'\n$insertPoint.InsertXML($xml)\n\n# Remove the standalone \"1 Introduction\" heading paragraph (and its\n# bookmark) entirely -- the heading text now lives, unformatted, at the\n# start of the paragraph above.\n$headingPara.Range.Delete()\n"}
